# Adding profits (M_PL) block to the income table: columns R:Y, mirroring
# the existing M_%cit (B:I) / M_ETR (J:Q) blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sub-header labels repeated for each of the three top-level blocks.
$subHeaders = @(
    "GFA - Sales",
    "GFA - Sales + Emp",
    "IMF - Sales",
    "IMF - Sales + Emp",
    "OECD (20%) - Sales",
    "OECD (20%) - Sales + Emp",
    "OECD - Sales",
    "OECD - Sales + Emp"
)

# New top-level header "M_PL" spans R1:Y1 (merged), matching the existing
# "M_%cit" (B1:I1) and "M_ETR" (J1:Q1) merged headers. Merge before applying
# the per-cell formatting below so the merged range doesn't lose its inner
# borders afterwards.
$ws.Range("R1:Y1").Merge()

# Apply the same formatting used by the existing header cells (bold, centered,
# top-aligned, thin box border) to every new header cell in row 1 and row 2,
# columns R (18) through Y (25).
for ($c = 18; $c -le 25; $c++) {
    $cell1 = $ws.Cells.Item(1, $c)
    $cell1.Font.Bold = $true
    $cell1.HorizontalAlignment = -4108
    $cell1.VerticalAlignment = -4160
    $cell1.Borders.Item(10).LineStyle = 1
    $cell1.Borders.Item(7).LineStyle = 1
    $cell1.Borders.Item(8).LineStyle = 1
    $cell1.Borders.Item(9).LineStyle = 1

    $cell2 = $ws.Cells.Item(2, $c)
    $cell2.Font.Bold = $true
    $cell2.HorizontalAlignment = -4108
    $cell2.VerticalAlignment = -4160
    $cell2.Borders.Item(10).LineStyle = 1
    $cell2.Borders.Item(7).LineStyle = 1
    $cell2.Borders.Item(8).LineStyle = 1
    $cell2.Borders.Item(9).LineStyle = 1
    $cell2.Value = $subHeaders[$c - 18]
}

$ws.Cells.Item(1, 18).Value = "M_PL"

# New data values for the M_PL block (rows 4-8, columns R-Y).
$plData = @{
    4 = @(957691691302, 957942582918, 955500195836, 955751087452, 1007534436142, 1007534436142, 1007534436142, 1007534436142)
    5 = @(-269766813, -269766813, -269766813, -269766813, -269766813, -269766813, -269766813, -269766813)
    6 = @(1507896554, 15012781826, 1074779159, 16007585528, 20228669958, 20228669958, 20228669958, 20228669958)
    7 = @(2770915272, -7824629507, -6943500091, -7824629507, -8693540732, -8693540732, -8693540732, -8693540732)
    8 = @(37772977295, 43238497013, 37492196667, 43740019703, 45733381438, 45733381438, 45733381438, 45733381438)
}

foreach ($r in $plData.Keys) {
    $vals = $plData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 18 + $i).Value = $vals[$i]
    }
}
